# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 16 - India
$ws.Range("B16").Value = 67724
$ws.Range("C16").Value = 563
$ws.Range("D16").Value = 21155
$ws.Range("E16").Value = 44354

# Row 19 - Paises Bajos
$ws.Range("B19").Value = 42788
$ws.Range("C19").Value = 161
$ws.Range("E19").Value = 37082
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = 5456

# Row 26 - Portugal
$ws.Range("B26").Value = 27679
$ws.Range("C26").Value = 98
$ws.Range("E26").Value = 23986
$ws.Range("G26").Value = 9
$ws.Range("H26").Value = 1144

# Row 27 - Suecia
$ws.Range("B27").Value = 26670
$ws.Range("C27").Value = 348
$ws.Range("E27").Value = 18443
$ws.Range("G27").Value = 31
$ws.Range("H27").Value = 3256

# Row 57 - Finlandia
$ws.Range("E57").Value = 1713
$ws.Range("F57").Value = 44
$ws.Range("G57").Value = 4
$ws.Range("H57").Value = 271

# Row 76 - Croacia
$ws.Range("B76").Value = 2196
$ws.Range("C76").Value = 9
$ws.Range("D76").Value = 1784
$ws.Range("E76").Value = 321
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 91
